# Apply updated hourly regression-with-fixed-effects coefficient table values
# (efficiencies using hourly regression with FE)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 0.09435480859084534
# Row 3
$ws.Range("B3").Value = 0.002403912873815977
$ws.Range("C3").Value = 0.0006046928049735363
$ws.Range("D3").Value = 4.149105628736279
$ws.Range("E3").Value = 0.06732253236063347
$ws.Range("F3").Value = 0.001218730312147687
$ws.Range("G3").Value = 0.00358909543548427
$ws.Range("H3").Value = 0.09675872146466132
# Row 4
$ws.Range("B4").Value = 0.01024091117930792
$ws.Range("C4").Value = 0.001011414943208974
$ws.Range("D4").Value = 9.051232571738625
$ws.Range("E4").Value = 0.04769876514732291
$ws.Range("F4").Value = 0.008258563553724391
$ws.Range("G4").Value = 0.01222325880489146
$ws.Range("H4").Value = 0.1045957197701533
# Row 5
$ws.Range("B5").Value = 0.06087690981765643
$ws.Range("C5").Value = 0.003646393340780202
$ws.Range("D5").Value = 13.0966092325068
$ws.Range("E5").Value = 0.000008297900597811221
$ws.Range("F5").Value = 0.05373007421576488
$ws.Range("G5").Value = 0.068023745419548
$ws.Range("H5").Value = 0.1552317184085018
# Row 6
$ws.Range("B6").Value = 0.148641638762907
$ws.Range("C6").Value = 0.00689419406925055
$ws.Range("D6").Value = 24.17799185899922
$ws.Range("E6").Value = 0.02453548909062416
$ws.Range("F6").Value = 0.1351292109030406
$ws.Range("G6").Value = 0.1621540666227735
$ws.Range("H6").Value = 0.2429964473537524
# Row 7
$ws.Range("B7").Value = 0.1556083542319699
$ws.Range("C7").Value = 0.009797452638638134
$ws.Range("D7").Value = 23.36333613834265
$ws.Range("E7").Value = 0.04681745132946349
$ws.Range("F7").Value = 0.136405632036056
$ws.Range("G7").Value = 0.1748110764278839
$ws.Range("H7").Value = 0.2499631628228153
# Row 8
$ws.Range("B8").Value = 0.1262119232944673
$ws.Range("C8").Value = 0.005053137174838146
$ws.Range("D8").Value = 26.80905749667661
$ws.Range("E8").Value = 0.03702528842349516
$ws.Range("F8").Value = 0.1163079221514899
$ws.Range("G8").Value = 0.1361159244374448
$ws.Range("H8").Value = 0.2205667318853127
# Row 9
$ws.Range("B9").Value = 0.1614645082822768
$ws.Range("C9").Value = 0.005072185175960149
$ws.Range("D9").Value = 30.07050980183883
$ws.Range("E9").Value = 0.01295126740128811
$ws.Range("F9").Value = 0.1515231627768194
$ws.Range("G9").Value = 0.1714058537877342
$ws.Range("H9").Value = 0.2558193168731221
# Row 10
$ws.Range("B10").Value = -0.09435480859084534
$ws.Range("C10").Value = 0.0004773167150480023
$ws.Range("D10").Value = -217.4873512448193
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = -0.0952903367084274
$ws.Range("G10").Value = -0.09341928047326332
# Row 11
$ws.Range("B11").Value = -0.04246719039222534
$ws.Range("C11").Value = 0.0005424364399391171
$ws.Range("D11").Value = -81.90375937554948
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = -0.04353035142435559
$ws.Range("G11").Value = -0.04140402936009509
$ws.Range("H11").Value = 0.05188761819862001
# Row 12
$ws.Range("B12").Value = -0.03550239386095641
$ws.Range("C12").Value = 0.0005317802979303487
$ws.Range("D12").Value = -68.34002413869803
$ws.Range("E12").Value = 0.0000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000001762030943082861
$ws.Range("F12").Value = -0.03654466915104662
$ws.Range("G12").Value = -0.03446011857086619
$ws.Range("H12").Value = 0.05885241472988893
# Row 13
$ws.Range("B13").Value = -0.02917258061474014
$ws.Range("C13").Value = 0.0005435910832097472
$ws.Range("D13").Value = -54.71148532831575
$ws.Range("E13").Value = 0.000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000000005107383341384822
$ws.Range("F13").Value = -0.03023800487301915
$ws.Range("G13").Value = -0.02810715635646114
$ws.Range("H13").Value = 0.0651822279761052
# Row 14
$ws.Range("B14").Value = -0.02555556142957584
$ws.Range("C14").Value = 0.0005324262492471209
$ws.Range("D14").Value = -47.46898931106173
$ws.Range("E14").Value = 0.0000000000000000000000000000000000000000000000000000000000001196211425284462
$ws.Range("F14").Value = -0.02659910287314212
$ws.Range("G14").Value = -0.02451201998600956
$ws.Range("H14").Value = 0.0687992471612695
# Row 15
$ws.Range("B15").Value = -0.0199283093706227
$ws.Range("C15").Value = 0.0005161643122518713
$ws.Range("D15").Value = -39.0977475129877
$ws.Range("E15").Value = 0.0000002316799870602134
$ws.Range("F15").Value = -0.02093997791512081
$ws.Range("G15").Value = -0.0189166408261246
$ws.Range("H15").Value = 0.07442649922022264
# Row 16
$ws.Range("B16").Value = -0.01753250648117852
$ws.Range("C16").Value = 0.0004802995299431071
$ws.Range("D16").Value = -36.67797658071094
$ws.Range("E16").Value = 0.0529967000641883
$ws.Range("F16").Value = -0.01847388104430327
$ws.Range("G16").Value = -0.01659113191805377
$ws.Range("H16").Value = 0.07682230210966683
# Row 17
$ws.Range("B17").Value = -0.01547955963961896
$ws.Range("C17").Value = 0.0004953222227543079
$ws.Range("D17").Value = -30.42279407529585
$ws.Range("E17").Value = 0.0000000000000000000001323335854612009
$ws.Range("F17").Value = -0.01645037819023642
$ws.Range("G17").Value = -0.0145087410890015
$ws.Range("H17").Value = 0.07887524895122638
# Row 18
$ws.Range("B18").Value = -0.01271862244494612
$ws.Range("C18").Value = 0.0005203507763497718
$ws.Range("D18").Value = -23.24264796306787
$ws.Range("E18").Value = 0.00005437121218992554
$ws.Range("F18").Value = -0.01373849633653106
$ws.Range("G18").Value = -0.01169874855336118
$ws.Range("H18").Value = 0.08163618614589922
# Row 19
$ws.Range("B19").Value = -0.0101841108751682
$ws.Range("C19").Value = 0.000537735591525884
$ws.Range("D19").Value = -17.67026099884727
$ws.Range("E19").Value = 0.07731248174601105
$ws.Range("F19").Value = -0.01123805888241611
$ws.Range("G19").Value = -0.009130162867920285
$ws.Range("H19").Value = 0.08417069771567715
# Row 20
$ws.Range("B20").Value = -0.008979698013693508
$ws.Range("C20").Value = 0.0005326090811532338
$ws.Range("D20").Value = -14.48067645808695
$ws.Range("E20").Value = 0.01922332833845312
$ws.Range("F20").Value = -0.01002359820127895
$ws.Range("G20").Value = -0.007935797826108065
$ws.Range("H20").Value = 0.08537511057715183
# Row 21
$ws.Range("B21").Value = -0.007800061518731668
$ws.Range("C21").Value = 0.0005297896184014143
$ws.Range("D21").Value = -12.32883015589398
$ws.Range("E21").Value = 0.1328930758081338
$ws.Range("F21").Value = -0.008838435593977047
$ws.Range("G21").Value = -0.006761687443486292
$ws.Range("H21").Value = 0.08655474707211368
# Row 22
$ws.Range("B22").Value = -0.00601869124245415
$ws.Range("C22").Value = 0.0005182453563018003
$ws.Range("D22").Value = -9.120369155499743
$ws.Range("E22").Value = 0.09894062606354796
$ws.Range("F22").Value = -0.007034438714827346
$ws.Range("G22").Value = -0.005002943770080956
$ws.Range("H22").Value = 0.0883361173483912
# Row 23
$ws.Range("B23").Value = -0.004897516296441459
$ws.Range("C23").Value = 0.0005058514108379643
$ws.Range("D23").Value = -7.579112915850029
$ws.Range("E23").Value = 0.08025350423650958
$ws.Range("F23").Value = -0.005888971969813844
$ws.Range("G23").Value = -0.003906060623069074
$ws.Range("H23").Value = 0.08945729229440388
# Row 24
$ws.Range("B24").Value = -0.004828631838060799
$ws.Range("C24").Value = 0.0004870030769324599
$ws.Range("D24").Value = -7.40979762896782
$ws.Range("E24").Value = 0.1496211128708858
$ws.Range("F24").Value = -0.005783145052994634
$ws.Range("G24").Value = -0.003874118623126968
$ws.Range("H24").Value = 0.08952617675278454
# Row 25
$ws.Range("B25").Value = -0.002795294661318059
$ws.Range("C25").Value = 0.0004643788831731986
$ws.Range("D25").Value = -5.19344164207013
$ws.Range("E25").Value = 0.1400155493671354
$ws.Range("F25").Value = -0.003705465207221732
$ws.Range("G25").Value = -0.001885124115414385
$ws.Range("H25").Value = 0.09155951392952728
# Row 26
$ws.Range("B26").Value = 0.2194895188580215
$ws.Range("C26").Value = 0.001200816808151769
$ws.Range("D26").Value = 182.7918824392778
$ws.Range("E26").Value = 0.00000000002442464564115682
$ws.Range("F26").Value = 0.2171359512597167
$ws.Range("G26").Value = 0.22184308645632633
$ws.Range("H26").Value = 0.31384432744886687
